$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 162.5768232345581

$ws.Range("A4").Value = 18814.8955
$ws.Range("B4").Value = 17451
$ws.Range("F4").Value = 6411.1545
$ws.Range("G4").Value = 6410
